$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.968.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.40%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.247.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.67%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'231.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.02%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.632"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.65%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'63.23"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.23%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.00%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.448"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +5.49%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0981"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.12%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'57.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.77%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'26.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.20%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.58%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.581.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.62%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'15.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.27%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'6.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.58%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.829"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.66%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.248.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.41%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'43.869.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.25%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0989"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +3.38%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'72.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.72%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.95%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'247.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.85%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +0.07%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -7.19%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +21.79%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.56%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.55%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'171.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.27%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'20.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.77%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.68%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -2.64%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +2.40%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.0685"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.87%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'4.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.17%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.96%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.32%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -2.68%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -3.72%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.96%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.11%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.000225"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.43%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'8.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.32%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'17.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.01%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'97.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.10%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -2.17%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -2.49%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'4.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -7.23%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.436.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.39%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -2.36%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +1.19%  "
$ws.Range("E51").Style = "Normal"
